$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Delete the whole "Vergonha e honra, Vestes sacerdotais, ..." list
#    paragraph (index-of-terms paragraph right after the "V" heading).
# ------------------------------------------------------------------
$listPara = $d.Paragraphs.Item(11)
if ($listPara.Range.Text -notmatch "Vergonha e honra, Vestes sacerdotais") {
    throw "Paragraph 11 did not match expected list paragraph text: $($listPara.Range.Text)"
}
$listPara.Range.Delete()

# ------------------------------------------------------------------
# 2. Delete the "This PDF version is provided under the same license."
#    paragraph.
# ------------------------------------------------------------------
$pdfPara = $d.Paragraphs.Item(6)
if ($pdfPara.Range.Text -notmatch "This PDF version is provided under the same license") {
    throw "Paragraph 6 did not match expected PDF paragraph text: $($pdfPara.Range.Text)"
}
$pdfPara.Range.Delete()

# ------------------------------------------------------------------
# 3. Rewrite the license/attribution paragraph (formerly about
#    "Termos Chave (Biblica)" + hyperlinks) with the new Biblica Study
#    Notes attribution text.
# ------------------------------------------------------------------
$licensePara = $d.Paragraphs.Item(5)
if ($licensePara.Range.Text -notmatch "is based on") {
    throw "Paragraph 5 did not match expected license paragraph text: $($licensePara.Range.Text)"
}
$fullRange = $licensePara.Range
$content = $d.Range($fullRange.Start, $fullRange.End - 1)
$content.Text = ""

$cursor = $d.Range($content.Start, $content.Start)
$cursor.Text = "Biblica Study Notes (Key Terms)"
$cursor.Font.Bold = 1

$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.Text = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. "
$cursor.Font.Bold = 0

$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.Text = "Biblica Study Notes"
$cursor.Font.Bold = 0

$cursor = $d.Range($cursor.End, $cursor.End)
$cursor.Text = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."
$cursor.Font.Bold = 0

# ------------------------------------------------------------------
# 4. Delete the "License Information" Heading2 paragraph.
# ------------------------------------------------------------------
$licenseHeadingPara = $d.Paragraphs.Item(4)
if ($licenseHeadingPara.Range.Text -notmatch "License Information") {
    throw "Paragraph 4 did not match expected heading text: $($licenseHeadingPara.Range.Text)"
}
$licenseHeadingPara.Range.Delete()

Write-Output "Done"
